$p = $ppt.ActivePresentation

# --- 1) Refresh the "datetimeFigureOut" date placeholder text on the slide
#        master and every slide layout: 2/17/2015 -> 2/21/2015 -------------
$master = $p.SlideMaster

$mshapes = $master.Shapes
for ($i = 1; $i -le $mshapes.Count; $i++) {
    $sh = $mshapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "2/17/2015") {
            $tr.Text = "2/21/2015"
        }
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    $lshapes = $lay.Shapes
    for ($i = 1; $i -le $lshapes.Count; $i++) {
        $sh = $lshapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2/17/2015") {
                $tr.Text = "2/21/2015"
            }
        }
    }
}

# --- 2) Slide 4: merge the " " and "+ " runs into a single " + " run in the
#        "li_M + Delta l_M" label (Rectangle 25) --------------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item("Rectangle 25")
$tr4 = $sh4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(1)
$sub4 = $para4.Characters(4, 3)
$sub4.Text = " + "

# --- 3) Slide 7: fix the node max-capacity numbers -------------------------
$s7 = $p.Slides.Item(7)

$sh7a = $s7.Shapes.Item(1)
$para7a = $sh7a.TextFrame.TextRange.Paragraphs(2)
$sub7a = $para7a.Characters($para7a.Length, 1)
$sub7a.Text = "8"

$sh7b = $s7.Shapes.Item(2)
$para7b = $sh7b.TextFrame.TextRange.Paragraphs(2)
$sub7b = $para7b.Characters($para7b.Length, 1)
$sub7b.Text = "5"

$sh7c = $s7.Shapes.Item(3)
$para7c = $sh7c.TextFrame.TextRange.Paragraphs(2)
$sub7c = $para7c.Characters($para7c.Length, 1)
$sub7c.Text = "8"

$sh7d = $s7.Shapes.Item(4)
$para7d = $sh7d.TextFrame.TextRange.Paragraphs(2)
$sub7d = $para7d.Characters($para7d.Length - 1, 2)
$sub7d.Text = "8"
